# Insert a new price-record row for Ajo (Chino / Primera) dated 2023-04-11
# (Excel serial 45027) right after the existing row 276 (which already
# holds the 2023-03-31 record). This mirrors the weekly price log pattern:
# a new week's row is added and all the historical rows below it shift
# down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 277..386 down to 278..387, leaving row 277 blank.
$ws.Rows(277).Insert()

# Populate the newly inserted row 277 with the new weekly record.
$ws.Range("A277").Value = 7
$ws.Range("B277").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C277").Value = 'Ñuble'
$ws.Range("D277").Value = 45027
$ws.Range("E277").Value = 16
$ws.Range("F277").Value = 100112003
$ws.Range("G277").Value = 'Ajo'
$ws.Range("H277").Value = 'Chino'
$ws.Range("I277").Value = 'Primera'
$ws.Range("J277").Value = 60
$ws.Range("K277").Value = 17000
$ws.Range("L277").Value = 18000
$ws.Range("M277").Value = 17500
$ws.Range("N277").Value = '$/caja 10 kilos'
$ws.Range("O277").Value = 'China'
$ws.Range("P277").Value = 1750
$ws.Range("Q277").Value = 10
$ws.Range("R277").Value = 'Hortaliza'
